# Update "doctor_MA" (column AF) result values on Sheet1, rows 4-13
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF4").Value  = 0.6899999999999999
$ws.Range("AF5").Value  = 0.96
$ws.Range("AF6").Value  = 0.803
$ws.Range("AF7").Value  = 0.89
$ws.Range("AF8").Value  = 0.866
$ws.Range("AF9").Value  = 0.72
$ws.Range("AF10").Value = 0.96
$ws.Range("AF11").Value = 0.96
$ws.Range("AF12").Value = 1.292
$ws.Range("AF13").Value = 1.68
